$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# Copy the header-row style (bold, bordered, centered/top) onto A2:A10
# so the new index column matches style index used by the existing header (s="1").
$ws.Cells.Item(1, 2).Copy()
$ws.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'Date Crown Lulu Dates'
$ws.Cells.Item(2, 3).Value = '400 gm'
$ws.Cells.Item(2, 4).Value = '৳ 249'
$ws.Cells.Item(2, 5).Value = '*Features an elongated oval *shape and smooth, glossy skin *Rich in fibers *High in potassium *Low in sodium *Source of magnesium *Naturally sweet flavour makes these dates a perfect choice for adding to cakes, cookies, and more *Ideal on-the-go snacking option *Can be stored in an airtight container to preserve freshness *Vegan *Gluten Free *No added sugar and preservatives '
$ws.Cells.Item(2, 6).Value = 'Ramadan'

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'Ajwa Premium Dates'
$ws.Cells.Item(3, 3).Value = '500 gm'
$ws.Cells.Item(3, 4).Value = '৳ 699'
$ws.Cells.Item(3, 5).Value = 'Origin: Saudi Arabia Dry,Semi-Moist,Oragnic,Dairy Free,Nut Free,Egg Free Ajwa dates Ingredients: Mabroom/Morium Dates/ Khejur Imported Dates. Halal Net Weight: 500gm Ajwa Dates (0.25 cup) contains 30g total carbs, 27g net carbs, 0g fat, 1g protein, and 120 calories'
$ws.Cells.Item(3, 6).Value = 'Ramadan'

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 'Crown Dabbas Dates'
$ws.Cells.Item(4, 3).Value = '1 kg'
$ws.Cells.Item(4, 4).Value = '৳ 549'
$ws.Cells.Item(4, 6).Value = 'Ramadan'

# Row 5
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'Chaldal Basic Dishwashing Bar'
$ws.Cells.Item(5, 3).Value = '100 gm'
$ws.Cells.Item(5, 4).Value = '৳ 12'
$ws.Cells.Item(5, 5).Value = 'Chaldal Basic Dishwashing Bar with the power of 100 lemons helps to clean tough grease the fastest. It gives you a pleasant cleaning experience with its refreshing lemon fragrance. It removes stains easily. It is tough on stains, gentle on your hands. Just take a little bit of Chaldal Basic Dishwashing Bar on your scrubber and clean your cooking vessels. So go ahead and enjoy the Chaldal Basic Dishwashing Bar experience. '
$ws.Cells.Item(5, 6).Value = 'Flash Sales'

# Row 6
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'Chaldal Basic Dishwashing Bar'
$ws.Cells.Item(6, 3).Value = '300 gm'
$ws.Cells.Item(6, 4).Value = '৳ 29'
$ws.Cells.Item(6, 5).Value = 'Chaldal Basic Dishwashing Bar with the power of 100 lemons helps to clean tough grease the fastest. It gives you a pleasant cleaning experience with its refreshing lemon fragrance. It removes stains easily. It is tough on stains, and gentle on your hands. Just take a little bit of Chaldal Basic Dishwashing Bar on your scrubber and clean your cooking vessels. So go ahead and enjoy the Chaldal Basic Dishwashing Bar experience.'
$ws.Cells.Item(6, 6).Value = 'Flash Sales'

# Row 7
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'Harpic Liquid Toilet Cleaner Original'
$ws.Cells.Item(7, 3).Value = '1 ltr'
$ws.Cells.Item(7, 4).Value = '৳ 169'
$ws.Cells.Item(7, 5).Value = 'Highlights: Removes tough stains Kills 99.9% of germs Leaves toilet sparkly Clean, hygienic and visibly neat Description: Harpic Toilet Cleaner Liquid 1L is a specialized all-in-one product, which is the ultimate one stop solution for all your toilet cleaning needs. Unlike ordinary toilet cleaners Harpic Power Plus has a thick liquid formula that clings to the toilet surface to clean deeply from the rim to the u-bend. Harpic toilet cleaner removes 99.9% of germs. Usage Information: Press Harpic Toilet Cleaner Liquid cap and twist to open, squeeze liquid around the toilet bowl and rim, leave for 5-10mins, finally scrub gently using a brush and flush. Safety Measure: Always use Harpic separately. Do not mix with other products. This product is not edible and harmful for eyes and skin. Wash with water immediately in case of contact with Skin or Eyes and consult a doctor. Read label for detail. About the Manufacturer: Harpic is a registered Trademark of Reckitt PLC headquartered in England. The first Harpic product was manufactured during 1920s in England. Harpic is the number 1 selling toilet cleaner in the world. Harpic products are made in Bangladesh by Reckitt (Bangladesh) PLC.'
$ws.Cells.Item(7, 6).Value = 'Flash Sales'

# Row 8
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 'Date Crown Lulu Dates'
$ws.Cells.Item(8, 3).Value = '400 gm'
$ws.Cells.Item(8, 4).Value = '৳ 249'
$ws.Cells.Item(8, 5).Value = '*Features an elongated oval *shape and smooth, glossy skin *Rich in fibers *High in potassium *Low in sodium *Source of magnesium *Naturally sweet flavour makes these dates a perfect choice for adding to cakes, cookies, and more *Ideal on-the-go snacking option *Can be stored in an airtight container to preserve freshness *Vegan *Gluten Free *No added sugar and preservatives '
$ws.Cells.Item(8, 6).Value = 'Popular'

# Row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'Ajwa Premium Dates'
$ws.Cells.Item(9, 3).Value = '500 gm'
$ws.Cells.Item(9, 4).Value = '৳ 699'
$ws.Cells.Item(9, 5).Value = 'Origin: Saudi Arabia Dry,Semi-Moist,Oragnic,Dairy Free,Nut Free,Egg Free Ajwa dates Ingredients: Mabroom/Morium Dates/ Khejur Imported Dates. Halal Net Weight: 500gm Ajwa Dates (0.25 cup) contains 30g total carbs, 27g net carbs, 0g fat, 1g protein, and 120 calories'
$ws.Cells.Item(9, 6).Value = 'Popular'

# Row 10
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 'Crown Dabbas Dates'
$ws.Cells.Item(10, 3).Value = '1 kg'
$ws.Cells.Item(10, 4).Value = '৳ 549'
$ws.Cells.Item(10, 6).Value = 'Popular'

Write-Output "done"
